$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.06056433333333333
$ws.Range("H2").Value = 0.181693
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6153926666666667
$ws.Range("N2").Value = 1.846178
$ws.Range("O2").Value = 0.172697186719763
$ws.Range("P2").Value = 0.172697186719763
$ws.Range("Q2").Value = 0.03727084659488889
$ws.Range("R2").Value = 0.335437619354
$ws.Range("S2").Value = 0.172697186719763
$ws.Range("T2").Value = 0.172697186719763

# Update row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 0.06056433333333333
$ws.Range("H3").Value = 0.181693
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.433374
$ws.Range("N3").Value = 1.300122
$ws.Range("O3").Value = 0.1216174235596306
$ws.Range("P3").Value = 0.1216174235596306
$ws.Range("Q3").Value = 0.026247007394
$ws.Range("R3").Value = 0.236223066546
$ws.Range("S3").Value = 0.1216174235596306
$ws.Range("T3").Value = 0.1216174235596306

# Update row 4
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 0.06056433333333333
$ws.Range("H4").Value = 0.181693
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.514653666666667
$ws.Range("N4").Value = 7.543961
$ws.Range("O4").Value = 0.7056853897206065
$ws.Range("P4").Value = 0.7056853897206065
$ws.Range("Q4").Value = 0.1522983228858889
$ws.Range("R4").Value = 1.370684905973
$ws.Range("S4").Value = 0.7056853897206065
$ws.Range("T4").Value = 0.7056853897206065

# Remove row 5 entirely (shifts rows up, deletes the ECs row data)
$ws.Range("A5:T5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
